# "Loan RBI, Variable Instalments"
#
# On the "Repayment schedule" sheet a new (blank) column is inserted at N,
# pushing the former N/O/P ("Late" / "Date" / "Outstanding") columns one
# place to the right (to O/P/Q). The new column inherits the column width
# of the column immediately to its left (M) and stays style-only / empty
# of values. The sheet also becomes the active sheet/tab, with cell K17
# selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (drives tabSelected on sheets + activeTab on
# the workbook's bookViews).
$ws.Activate()

# Insert a new blank column before the existing column N ("Late"),
# shifting N:P to O:Q.
$beforeWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert() | Out-Null

# The newly inserted column picks up the same width as its left neighbour.
$ws.Columns("N:N").ColumnWidth = $beforeWidth

# Update the selected cell on the sheet.
$ws.Range("K17").Select() | Out-Null
